$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the refreshed cryptocurrency market data.
# D-column (Price) values are forced to Text so Excel does not
# auto-convert numeric-looking strings (e.g. "310.49") into numbers,
# matching the original inline-string storage; formatting is cleared
# afterward so no stray style is left on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.086.59'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.305.18'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.87%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.49'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.35'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +6.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.537'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.11%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.524'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +6.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.98'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0818'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.49%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.16'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +7.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.660.40'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.99'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.301.00'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.809'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.038.53'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.52'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0924'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.09'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.51'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.55'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('E25').Value = '  +3.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.67'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.36%  '
$ws.Range('E28').Value = '  +10.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.67'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.63'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.44'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.31'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.69'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.95%  '
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.107'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('E39').Value = '  +1.70%  '
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.24'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.83%  '
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.978.57'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0289'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.13'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.03'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.30%  '
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.96'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +18.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.62'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.528.28'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.53'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.08%  '
